# Auto-generated script applying scheduled market-data refresh to Aegis_Profits (Leve Profit) sheets.
# For each sheet, update computed market-price / profit columns (H-N) for specific Leve rows
# to match the latest Universalis market snapshot pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 812.2
$ws.Range("I19").Value = 615.0833
$ws.Range("J19").Value = 943.6111
$ws.Range("K19").Value = 615.0833
$ws.Range("L19").Value = 943.6111
$ws.Range("M19").Value = -440.0833
$ws.Range("N19").Value = -1293.6111
$ws.Range("H40").Value = 1964.9062
$ws.Range("I40").Value = 1895.2941
$ws.Range("J40").Value = 2043.8
$ws.Range("K40").Value = 1895.2941
$ws.Range("L40").Value = 2043.8
$ws.Range("M40").Value = -1720.2941
$ws.Range("N40").Value = -2393.8
$ws.Range("H108").Value = 21374.5
$ws.Range("J108").Value = 21374.5
$ws.Range("L108").Value = 21374.5
$ws.Range("N108").Value = -29054.5
$ws.Range("H113").Value = 60192.883
$ws.Range("J113").Value = 1474.1666
$ws.Range("L113").Value = 1474.1666
$ws.Range("N113").Value = -7982.1666
$ws.Range("H125").Value = 2389.5
$ws.Range("I125").Value = 2220.75
$ws.Range("J125").Value = 2727
$ws.Range("K125").Value = 19986.75
$ws.Range("L125").Value = 24543
$ws.Range("M125").Value = -17526.75
$ws.Range("N125").Value = -29463
$ws.Range("H135").Value = 947.7308
$ws.Range("I135").Value = 814.8261
$ws.Range("K135").Value = 7333.4349
$ws.Range("M135").Value = -4798.4349
$ws.Range("H137").Value = 1466.8636
$ws.Range("I137").Value = 1412.0667
$ws.Range("J137").Value = 1584.2858
$ws.Range("K137").Value = 4236.2001
$ws.Range("L137").Value = 4752.857400000001
$ws.Range("M137").Value = -1686.2001
$ws.Range("N137").Value = -9852.857400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26516.174
$ws.Range("I32").Value = 4447.4917
$ws.Range("J32").Value = 156721.4
$ws.Range("K32").Value = 4447.4917
$ws.Range("L32").Value = 156721.4
$ws.Range("M32").Value = -4160.4917
$ws.Range("N32").Value = -157295.4
$ws.Range("H61").Value = 1573.9688
$ws.Range("I61").Value = 1479.5
$ws.Range("J61").Value = 1983.3334
$ws.Range("K61").Value = 1479.5
$ws.Range("L61").Value = 1983.3334
$ws.Range("M61").Value = -1267.5
$ws.Range("N61").Value = -2407.3334
$ws.Range("H97").Value = 42985.875
$ws.Range("I97").Value = 44371.74
$ws.Range("K97").Value = 44371.74
$ws.Range("M97").Value = -43875.74
$ws.Range("H98").Value = 11000
$ws.Range("J98").Value = 11000
$ws.Range("L98").Value = 11000
$ws.Range("N98").Value = -16990
$ws.Range("H130").Value = 19000
$ws.Range("J130").Value = 19000
$ws.Range("L130").Value = 19000
$ws.Range("N130").Value = -29040
$ws.Range("H132").Value = 2332.653
$ws.Range("I132").Value = 2164.8438
$ws.Range("J132").Value = 2648.5293
$ws.Range("K132").Value = 6494.5314
$ws.Range("L132").Value = 7945.5879
$ws.Range("M132").Value = -3964.5314
$ws.Range("N132").Value = -13005.5879
$ws.Range("H136").Value = 1573.9688
$ws.Range("I136").Value = 1479.5
$ws.Range("J136").Value = 1983.3334
$ws.Range("K136").Value = 4438.5
$ws.Range("L136").Value = 5950.0002
$ws.Range("M136").Value = -1888.5
$ws.Range("N136").Value = -11050.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 330.44446
$ws.Range("I22").Value = 309.25
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 309.25
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -136.25
$ws.Range("N22").Value = -846
$ws.Range("H94").Value = 742.5
$ws.Range("I94").Value = 485.66666
$ws.Range("K94").Value = 485.66666
$ws.Range("M94").Value = -34.66665999999998
$ws.Range("H105").Value = 183547.36
$ws.Range("I105").Value = 168635
$ws.Range("J105").Value = 201442.2
$ws.Range("K105").Value = 168635
$ws.Range("L105").Value = 201442.2
$ws.Range("M105").Value = -166888
$ws.Range("N105").Value = -204936.2
$ws.Range("H107").Value = 142924460
$ws.Range("I107").Value = 250114540
$ws.Range("J107").Value = 4337.3335
$ws.Range("K107").Value = 250114540
$ws.Range("L107").Value = 4337.3335
$ws.Range("M107").Value = -250112620
$ws.Range("N107").Value = -8177.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 54.666668
$ws.Range("I2").Value = 54.666668
$ws.Range("K2").Value = 54.666668
$ws.Range("M2").Value = 58.333332
$ws.Range("H99").Value = 30748.5
$ws.Range("I99").Value = 8990
$ws.Range("K99").Value = 8990
$ws.Range("M99").Value = -7492
$ws.Range("H125").Value = 42500
$ws.Range("J125").Value = 42500
$ws.Range("L125").Value = 42500
$ws.Range("N125").Value = -47420
$ws.Range("H126").Value = 30748.5
$ws.Range("I126").Value = 8990
$ws.Range("K126").Value = 26970
$ws.Range("M126").Value = -24500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 337160
$ws.Range("I120").Value = 337160
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 1011480
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -1006642
$ws.Range("N120").ClearContents()
$ws.Range("H131").Value = 811.84
$ws.Range("I131").Value = 481.15384
$ws.Range("J131").Value = 861.25287
$ws.Range("K131").Value = 1443.46152
$ws.Range("L131").Value = 2583.75861
$ws.Range("M131").Value = 3596.53848
$ws.Range("N131").Value = -12663.75861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 75.947365
$ws.Range("I2").Value = 39.5
$ws.Range("J2").Value = 138.42857
$ws.Range("K2").Value = 39.5
$ws.Range("L2").Value = 138.42857
$ws.Range("M2").Value = 73.5
$ws.Range("N2").Value = -364.42857
$ws.Range("H70").Value = 123737.3
$ws.Range("I70").Value = 226778.33
$ws.Range("J70").Value = 7816.125
$ws.Range("K70").Value = 226778.33
$ws.Range("L70").Value = 7816.125
$ws.Range("M70").Value = -226508.33
$ws.Range("N70").Value = -8356.125
$ws.Range("H73").Value = 123737.3
$ws.Range("I73").Value = 226778.33
$ws.Range("J73").Value = 7816.125
$ws.Range("K73").Value = 226778.33
$ws.Range("L73").Value = 7816.125
$ws.Range("M73").Value = -225842.33
$ws.Range("N73").Value = -9688.125
$ws.Range("H80").Value = 166668860
$ws.Range("I80").Value = 333335400
$ws.Range("K80").Value = 333335400
$ws.Range("M80").Value = -333334402
$ws.Range("H83").Value = 166668860
$ws.Range("I83").Value = 333335400
$ws.Range("K83").Value = 1666677000
$ws.Range("M83").Value = -1666672008
$ws.Range("H132").Value = 3856.0625
$ws.Range("I132").Value = 3462.5
$ws.Range("J132").Value = 4249.625
$ws.Range("K132").Value = 10387.5
$ws.Range("L132").Value = 12748.875
$ws.Range("M132").Value = -7857.5
$ws.Range("N132").Value = -17808.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 4999.5
$ws.Range("J10").Value = 4999.5
$ws.Range("L10").Value = 4999.5
$ws.Range("N10").Value = -5279.5
$ws.Range("H22").Value = 1522.3889
$ws.Range("I22").Value = 4066.3333
$ws.Range("J22").Value = 1013.6
$ws.Range("K22").Value = 4066.3333
$ws.Range("L22").Value = 1013.6
$ws.Range("M22").Value = -3771.3333
$ws.Range("N22").Value = -1603.6
$ws.Range("H27").Value = 1522.3889
$ws.Range("I27").Value = 4066.3333
$ws.Range("J27").Value = 1013.6
$ws.Range("K27").Value = 4066.3333
$ws.Range("L27").Value = 1013.6
$ws.Range("M27").Value = -3959.3333
$ws.Range("N27").Value = -1227.6
$ws.Range("H68").Value = 3814.875
$ws.Range("I68").Value = 1645.5714
$ws.Range("J68").Value = 19000
$ws.Range("K68").Value = 1645.5714
$ws.Range("L68").Value = 19000
$ws.Range("M68").Value = -896.5714
$ws.Range("N68").Value = -20498
$ws.Range("H71").Value = 3814.875
$ws.Range("I71").Value = 1645.5714
$ws.Range("J71").Value = 19000
$ws.Range("K71").Value = 8227.857
$ws.Range("L71").Value = 95000
$ws.Range("M71").Value = -4483.857
$ws.Range("N71").Value = -102488
$ws.Range("H81").Value = 36665.5
$ws.Range("J81").Value = 36665.5
$ws.Range("L81").Value = 36665.5
$ws.Range("N81").Value = -38661.5
$ws.Range("H84").Value = 36665.5
$ws.Range("J84").Value = 36665.5
$ws.Range("L84").Value = 109996.5
$ws.Range("N84").Value = -119980.5
$ws.Range("H95").Value = 13600
$ws.Range("J95").Value = 13600
$ws.Range("L95").Value = 13600
$ws.Range("N95").Value = -19092
$ws.Range("H128").Value = 38462.5
$ws.Range("J128").Value = 38462.5
$ws.Range("L128").Value = 38462.5
$ws.Range("N128").Value = -48422.5
$ws.Range("H132").Value = 4234.0527
$ws.Range("I132").Value = 4788.4287
$ws.Range("J132").Value = 2681.8
$ws.Range("K132").Value = 14365.2861
$ws.Range("L132").Value = 8045.400000000001
$ws.Range("M132").Value = -11835.2861
$ws.Range("N132").Value = -13105.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 27298.334
$ws.Range("J117").Value = 27298.334
$ws.Range("L117").Value = 27298.334
$ws.Range("N117").Value = -36476.334
$ws.Range("H122").Value = 1843.1111
$ws.Range("I122").Value = 1272.5714
$ws.Range("J122").Value = 2206.182
$ws.Range("K122").Value = 3817.7142
$ws.Range("L122").Value = 6618.545999999999
$ws.Range("M122").Value = -1367.7142
$ws.Range("N122").Value = -11518.546
$ws.Range("H124").Value = 35995
$ws.Range("J124").Value = 35995
$ws.Range("L124").Value = 35995
$ws.Range("N124").Value = -45815
$ws.Range("H126").Value = 1909.3334
$ws.Range("I126").Value = 2496.8
$ws.Range("J126").Value = 1175
$ws.Range("K126").Value = 7490.400000000001
$ws.Range("L126").Value = 3525
$ws.Range("M126").Value = -5020.400000000001
$ws.Range("N126").Value = -8465
